# Accountability Sheet - add "Sources" column (G) with reference links,
# highlight row 4 (task #2) to match row 3's "completed" style, and
# update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1: "Sources" (bold, centered - matches header style)
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Value = "Sources"

# G4 must be written before G3 so the new shared-strings are appended in
# the same order as the target workbook (cY4HiiFHO1o before bSaBmXFym30).
$ws.Range("G4").Interior.Color = 65535
$ws.Range("G4").HorizontalAlignment = -4108
$ws.Range("G4").Value = "https://www.youtube.com/watch?v=cY4HiiFHO1o"

$ws.Range("G3").Interior.Color = 65535
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").Value = "https://www.youtube.com/watch?v=bSaBmXFym30"

# Row 4 (task #2) gets the same yellow highlight formatting as row 3,
# including the centered date number format in column D.
$ws.Range("A4:F4").Interior.Color = 65535
$ws.Range("A4:F4").HorizontalAlignment = -4108
$ws.Range("D4").NumberFormat = "d-mmm-yy"

# Widen the new Sources column.
$ws.Columns.Item(7).ColumnWidth = 32

# Update the active cell selection.
$ws.Range("D11").Select()
